# "Cap nhat toi uu" — refresh the RSI/Khoi luong dashboard with a new data
# block (columns F:I, mirroring A:D) and re-point the comparison formulas
# in G39:G42 at the new block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New header / label block in F1:I9, mirroring A1:D9
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "5m"
$ws.Range("G1").Value = "30m"
$ws.Range("H1").Value = "1H"
$ws.Range("I1").Value = "4H"

$ws.Range("F2").Value = "Khoi luong"
$ws.Range("G2").Value = "Khoi luong"
$ws.Range("H2").Value = "Khoi luong"
$ws.Range("I2").Value = "Khoi luong"

$ws.Range("F5").Value = "RSI"
$ws.Range("G5").Value = "RSI"
$ws.Range("H5").Value = "RSI"
$ws.Range("I5").Value = "RSI"

$ws.Range("F8").Value = [char]0x2205
$ws.Range("G8").Value = [char]0x2205
$ws.Range("H8").Value = [char]0x2205
$ws.Range("I8").Value = [char]0x2205
$ws.Range("F9").Value = [char]0x2205
$ws.Range("G9").Value = [char]0x2205
$ws.Range("H9").Value = [char]0x2205
$ws.Range("I9").Value = [char]0x2205

# ---------------------------------------------------------------------
# Updated "Khoi luong" (volume) figures, A3:D4, plus new F3:I4 block
# ---------------------------------------------------------------------
$ws.Range("A3").Value = 58
$ws.Range("B3").Value = 425
$ws.Range("C3").Value = 425
$ws.Range("D3").Value = 5792

$ws.Range("A4").Value = 454.6
$ws.Range("B4").Value = 2510.1
$ws.Range("C4").Value = 4130.3500000000004
$ws.Range("D4").Value = 18167.400000000001

$ws.Range("F3").Value = 410
$ws.Range("G3").Value = 1110
$ws.Range("H3").Value = 1110
$ws.Range("I3").Value = 1110

$ws.Range("F4").Value = 549.9
$ws.Range("G4").Value = 2572.85
$ws.Range("H4").Value = 4315.45
$ws.Range("I4").Value = 17938.8

# ---------------------------------------------------------------------
# Updated RSI figures, A6:D7, plus new F6:I7 block
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 44.81
$ws.Range("B6").Value = 46.02
$ws.Range("C6").Value = 53.59
$ws.Range("D6").Value = 64.08

$ws.Range("A7").Value = 42.67
$ws.Range("B7").Value = 55.8
$ws.Range("C7").Value = 65.14
$ws.Range("D7").Value = 64.599999999999994

$ws.Range("F6").Value = 44.18
$ws.Range("G6").Value = 47.87
$ws.Range("H6").Value = 55.54
$ws.Range("I6").Value = 65.19

$ws.Range("F7").Value = 43.38
$ws.Range("G7").Value = 60.11
$ws.Range("H7").Value = 65.819999999999993
$ws.Range("I7").Value = 64.680000000000007

# ---------------------------------------------------------------------
# G39:G40 become formulas (were static constants 200 / 2000); G41:G42
# are brand-new formulas comparing the 1H / 4H "before" vs. "after" volume.
# ---------------------------------------------------------------------
$ws.Range("G39").Formula = "=A4-F4"
$ws.Range("G40").Formula = "=B4-G4"
$ws.Range("G41").Formula = "=C4-H4"
$ws.Range("G42").Formula = "=D4-I4"

# ---------------------------------------------------------------------
# Column widths for the new F:I columns
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 9.33
$ws.Columns.Item(7).ColumnWidth = 9.67
$ws.Columns.Item(8).ColumnWidth = 10.17
$ws.Columns.Item(9).ColumnWidth = 10.33

# ---------------------------------------------------------------------
# View state — scroll so row 29 is at the top and I38 is selected
# ---------------------------------------------------------------------
$ws.Range("I38").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
